$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A77").Value = "TM4C123gh6PM Special Pins."

$ws.Range("B78").Value = "PIN No"
$ws.Range("C78").Value = "Special Comments"

$ws.Range("A79").Value = "ADC"
$ws.Range("B79").Value = "PE 0"

$ws.Range("B80").Value = "PE1"

$ws.Range("B81").Value = "PE2"

$ws.Range("B82").Value = "PE3"

$ws.Range("B83").Value = "PD0"

$ws.Range("B84").Value = "PD1"

$ws.Range("B85").Value = "PD2"

$ws.Range("B86").Value = "PD3"

$ws.Range("B87").Value = "PE4"

$ws.Range("B88").Value = "PE5"

$ws.Range("B89").Value = "PB4"

$ws.Range("B90").Value = "PB5"

$ws.Range("A91").Value = "UART"
$ws.Range("B91").Value = "PA0 RX0"
$ws.Range("C91").Value = "PROGRAMMER UART"

$ws.Range("B92").Value = "PA1 TX0"
$ws.Range("C92").Value = "PROGRAMMER UART"

$ws.Range("B93").Value = "PC4,PB0 RX1"

$ws.Range("B94").Value = "PC5,PB1 TX1"

$ws.Range("B95").Value = "PD6 RX2"

$ws.Range("B96").Value = "PD7 TX2"

$ws.Range("B97").Value = "PC6 RX3"

$ws.Range("B98").Value = "PC7 RX3"

$ws.Range("B99").Value = "PC4 RX"

$ws.Range("B100").Value = "PC5 TX4"

$ws.Range("B101").Value = "PE4 RX5"

$ws.Range("B102").Value = "PE5 TX5"

$ws.Range("B103").Value = "PD4 RX6"

$ws.Range("B104").Value = "PD5 TX6"

$ws.Range("B105").Value = "PE0 RX7"

$ws.Range("B106").Value = "PE1 RX7"

$ws.Range("A107").Value = "PWM"
$ws.Range("B107").Value = "PIN Name"
$ws.Range("C107").Value = "Pin"
$ws.Range("D107").Value = "Generator"

$ws.Range("B108").Value = "M0PWM0"
$ws.Range("C108").Value = "PB6"
$ws.Range("D108").Value = 0

$ws.Range("B109").Value = "M0PWM1"
$ws.Range("C109").Value = "PB7"
$ws.Range("D109").Value = 0

$ws.Range("B110").Value = "M0PWM2"
$ws.Range("C110").Value = "PB4"
$ws.Range("D110").Value = 1

$ws.Range("B112").Value = "M0PWM3"
$ws.Range("C112").Value = "PB5"
$ws.Range("D112").Value = 1

$ws.Range("B113").Value = "M0PWM4"
$ws.Range("C113").Value = "PE4"
$ws.Range("D113").Value = 2

$ws.Range("B114").Value = "MOPWM5"
$ws.Range("C114").Value = "PE5"
$ws.Range("D114").Value = 2

$ws.Range("B115").Value = "M0PWM6"
$ws.Range("C115").Value = "PC4,PD0"
$ws.Range("D115").Value = 3

$ws.Range("B116").Value = "M0PWM7"
$ws.Range("C116").Value = "PC5,PD1"
$ws.Range("D116").Value = 3

$ws.Range("B117").Value = "M1PWM0"
$ws.Range("C117").Value = "PD0"
$ws.Range("D117").Value = 0

$ws.Range("B118").Value = "M1PWM1"
$ws.Range("C118").Value = "PD1"
$ws.Range("D118").Value = 0

$ws.Range("B119").Value = "M1PWM2"
$ws.Range("C119").Value = "PA6 PE4"
$ws.Range("D119").Value = 1

$ws.Range("B120").Value = "M1PWM3"
$ws.Range("C120").Value = "PA7 PE5"
$ws.Range("D120").Value = 1

$ws.Range("B121").Value = "M1PWM4"
$ws.Range("C121").Value = "PF0"
$ws.Range("D121").Value = 2

$ws.Range("B122").Value = "M1PWM5"
$ws.Range("C122").Value = "PF1"
$ws.Range("D122").Value = 2

$ws.Range("B123").Value = "M1PWM6"
$ws.Range("C123").Value = "PF2"
$ws.Range("D123").Value = 3

$ws.Range("B124").Value = "M1PWM7"
$ws.Range("C124").Value = "PF3"
$ws.Range("D124").Value = 3

$ws.Range("C73").Select()